$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text / naturally-non-numeric cell updates (names, links, volume %, multi-dot prices) ---
$ws.Range("D2").Value = "61.082.92"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.370.56"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("D9").Value = "2.370.47"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "2.782.87"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "61.021.46"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.379.29"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("E25").Value = "  -15.27%  "
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "2.472.35"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0883"
$ws.Range("E30").Value = "  -6.95%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E31").Value = "  -4.62%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  -8.17%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  -4.13%  "
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("E42").Value = "  +4.12%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  -8.25%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -5.22%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -0.29%  "

# --- Price cells whose new value would be auto-parsed as a number by Excel; force text type ---
# Step 1: mark as text format so the assigned value is stored as a string, not a number
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Step 2: assign the values (now interpreted as text because of the "@" format)
$ws.Range("D5").Value = "546.96"
$ws.Range("D6").Value = "137.84"
$ws.Range("D8").Value = "0.525"
$ws.Range("D10").Value = "0.107"
$ws.Range("D12").Value = "5.36"
$ws.Range("D13").Value = "0.347"
$ws.Range("D14").Value = "25.08"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("D19").Value = "10.77"
$ws.Range("D20").Value = "4.14"
$ws.Range("D21").Value = "320.37"
$ws.Range("D24").Value = "64.26"
$ws.Range("D25").Value = "1.65"
$ws.Range("D26").Value = "8.23"
$ws.Range("D29").Value = "8.11"
$ws.Range("D31").Value = "1.39"
$ws.Range("D33").Value = "498.62"
$ws.Range("D36").Value = "0.999"
$ws.Range("D37").Value = "4.68"
$ws.Range("D40").Value = "5.37"
$ws.Range("D41").Value = "18.56"
$ws.Range("D42").Value = "144.94"
$ws.Range("D44").Value = "41.15"
$ws.Range("D45").Value = "144.80"
$ws.Range("D46").Value = "3.59"
$ws.Range("D47").Value = "2.07"
$ws.Range("D49").Value = "19.24"
$ws.Range("D51").Value = "0.0909"

# Step 3: restore the original (General) number format on those cells
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "General"

Write-Host "Updated cryptos list"
